# Issue 205 and reorder slots in OrderedListItem

$wb = $excel.ActiveWorkbook

# Rename sheets AnalysisCategorization -> AnalysisOutputCategorization
# and AnalysisCategory -> AnalysisOutputCategory
$wsCategorization = $wb.Worksheets.Item("AnalysisCategorization")
$wsCategorization.Name = "AnalysisOutputCategorization"

$wsCategory = $wb.Worksheets.Item("AnalysisCategory")
$wsCategory.Name = "AnalysisOutputCategory"

# Update the ReportingEvent header label that referenced the old sheet/slot name
$wsReportingEvent = $wb.Worksheets.Item("ReportingEvent")
$wsReportingEvent.Range("G1").Value = "analysisOutputCategorizations"

# Reorder slots in OrderedListItem: sublist, analysisId, outputId -> analysisId, outputId, sublist
$wsOrderedListItem = $wb.Worksheets.Item("OrderedListItem")
$wsOrderedListItem.Range("C1").Value = "analysisId"
$wsOrderedListItem.Range("D1").Value = "outputId"
$wsOrderedListItem.Range("E1").Value = "sublist"
